# "Add files via upload"
#
# The Title paragraph currently reads:
#   "Test Case Document - Password Policy Including Special Character"
# This change prepends the word "Trial " (note the trailing space) in
# front of the existing title text, as its own run, so the title becomes:
#   "Trial Test Case Document - Password Policy Including Special Character"

$d = $word.ActiveDocument
$insertText = "Trial "

# The title is the very first paragraph in the document (style "Title").
$titleRange = $d.Paragraphs(1).Range
$titleRange.Collapse(1)             # wdCollapseStart -> start of the title
$titleRange.InsertBefore($insertText) # new text, inserted ahead of the old run

# Nudge the newly inserted text's direct character formatting (toggle a
# property on and back off) so the host keeps "Trial " as its own distinct
# run -- mirroring the original edit, which added a separate <w:r> for the
# inserted text -- instead of silently folding it into the pre-existing
# title run.
$newRun = $d.Range(0, $insertText.Length)
$newRun.Font.Bold = $true
$newRun.Font.Bold = $false
